# Refresh the "cryptos" price/volume snapshot (Price = column D, Volume(1h) = column E)
# for rows 2-51, matching the periodic GitHub Actions data-refresh commit.
# Price cells whose new text would otherwise be auto-parsed by Excel as a
# number (e.g. "241.14") are forced back to Text format first so they keep
# matching the sheet's existing plain-text price formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.920.49"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "2.218.19"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.14"
$ws.Range("E5").Value = "  -2.36%  "
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.50"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.67"
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0953"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.10"
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "2.549.08"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.29"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.838"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").Value = "2.227.84"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").Value = "41.865.11"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.22"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.15"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.69"
$ws.Range("E22").Value = "  +20.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.60"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  -6.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.73"
$ws.Range("E25").Value = "  +2.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.28"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("E32").Value = "  +8.11%  "
$ws.Range("E33").Value = "  -3.59%  "
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.23"
$ws.Range("E35").Value = "  -5.92%  "
$ws.Range("E36").Value = "  -10.74%  "
$ws.Range("E37").Value = "  -5.27%  "
$ws.Range("E38").Value = "  -4.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.62"
$ws.Range("E39").Value = "  -3.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "65.84"
$ws.Range("E40").Value = "  +5.95%  "
$ws.Range("E41").Value = "  -2.97%  "
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.198"
$ws.Range("E43").Value = "  -3.42%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.13"
$ws.Range("E45").Value = "  -2.89%  "
$ws.Range("E46").Value = "  -2.50%  "
$ws.Range("E47").Value = "  +3.23%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.17"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").Value = "2.423.68"
$ws.Range("E51").Value = "  -1.53%  "
